# landing page + slack
#
# 1) Bump the cached "datetimeFigureOut" footer field from 3/14/21 to
#    3/15/21 everywhere it is defined: the slide master, all eleven
#    slide layouts, and the notes master.
# 2) Shorten the "Releases and experiments powered by AIOps" callout on
#    slide 15 to "Releases and experiments powered by AI".

$p = $ppt.ActivePresentation

$oldDate = "3/14/21"
$newDate = "3/15/21"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
                $tr = $shape.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes

# Slide 15: trim "AIOps" down to "AI" in the rounded-rectangle callout
$slide15 = $p.Slides.Item(15)
for ($i = 1; $i -le $slide15.Shapes.Count; $i++) {
    $shape = $slide15.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "Releases and experiments powered by AIOps") {
            $tr.Text = "Releases and experiments powered by AI"
        }
    }
}
